$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that need to be forced to Text so Excel does not
# auto-convert numeric-looking strings (losing trailing zeros / formatting).
$textCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D36", "D38", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.626.32"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.526.60"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "314.13"
$ws.Range("E5").Value = "  +2.93%  "
$ws.Range("D6").Value = "95.43"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D10").Value = "36.47"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "2.919.45"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "15.71"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "2.539.15"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "0.861"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "42.680.50"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  -4.26%  "
$ws.Range("D20").Value = "6.69"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").Value = "0.0₃0971"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "71.23"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "254.54"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("D26").Value = "27.47"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "2.36"
$ws.Range("E28").Value = "  +12.97%  "
$ws.Range("D29").Value = "39.45"
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").Value = "5.91"
$ws.Range("E31").Value = "  -5.29%  "
$ws.Range("D32").Value = "155.66"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "20.17"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").Value = "0.0790"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "25.36"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "2.26"
$ws.Range("E41").Value = "  +9.59%  "
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "3.85"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0302"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "2.044.70"
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").Value = "85.54"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").Value = "8.90"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").Value = "74.75"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "2.778.49"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("E51").Value = "  -1.16%  "
